$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows before row 141, pushing existing rows 141-153 down to 143-155.
$ws.Rows.Item(141).Resize(2).Insert()

# Copy the date style (style index 2, format YYYY-MM-DD HH:MM:SS) used in column D
# for the newly inserted rows, matching the style of the surrounding rows.
$ws.Range("D143").Copy() | Out-Null
$ws.Range("D141:D142").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 141: new data record
$ws.Cells.Item(141, 1).Value = 10
$ws.Cells.Item(141, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(141, 3).Value = "La Araucanía"
$ws.Cells.Item(141, 4).Value = 44461
$ws.Cells.Item(141, 5).Value = 9
$ws.Cells.Item(141, 6).Value = 100112043
$ws.Cells.Item(141, 7).Value = "Pepino dulce"
$ws.Cells.Item(141, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(141, 9).Value = "Primera"
$ws.Cells.Item(141, 10).Value = 100
$ws.Cells.Item(141, 11).Value = 24000
$ws.Cells.Item(141, 12).Value = 24000
$ws.Cells.Item(141, 13).Value = 24000
$ws.Cells.Item(141, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(141, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(141, 16).Value = 1333
$ws.Cells.Item(141, 17).Value = 18
$ws.Cells.Item(141, 18).Value = "Hortaliza"

# Row 142: new data record
$ws.Cells.Item(142, 1).Value = 10
$ws.Cells.Item(142, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(142, 3).Value = "La Araucanía"
$ws.Cells.Item(142, 4).Value = 44461
$ws.Cells.Item(142, 5).Value = 9
$ws.Cells.Item(142, 6).Value = 100112043
$ws.Cells.Item(142, 7).Value = "Pepino dulce"
$ws.Cells.Item(142, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(142, 9).Value = "Segunda"
$ws.Cells.Item(142, 10).Value = 40
$ws.Cells.Item(142, 11).Value = 20000
$ws.Cells.Item(142, 12).Value = 20000
$ws.Cells.Item(142, 13).Value = 20000
$ws.Cells.Item(142, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(142, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(142, 16).Value = 1111
$ws.Cells.Item(142, 17).Value = 18
$ws.Cells.Item(142, 18).Value = "Hortaliza"
